$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-11-04 22:38:48"
